$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 108
$ws.Range("H108").Value = 31305
$ws.Range("J108").Value = 31305
$ws.Range("L108").Value = 31305
$ws.Range("N108").Value = -38985

# Row 109
$ws.Range("H109").Value = 33643.332
$ws.Range("J109").Value = 33643.332
$ws.Range("L109").Value = 33643.332
$ws.Range("N109").Value = -36417.332

# Row 114
$ws.Range("H114").Value = 35408
$ws.Range("J114").Value = 35408
$ws.Range("L114").Value = 35408
$ws.Range("N114").Value = -44086

# Row 120
$ws.Range("H120").Value = 49726
$ws.Range("J120").Value = 49726
$ws.Range("L120").Value = 49726
$ws.Range("N120").Value = -59402

# Row 124
$ws.Range("H124").Value = 40239.25
$ws.Range("J124").Value = 40239.25
$ws.Range("L124").Value = 40239.25
$ws.Range("N124").Value = -50059.25

# Row 130
$ws.Range("H130").Value = 52996
$ws.Range("J130").Value = 52996
$ws.Range("L130").Value = 52996
$ws.Range("N130").Value = -63036


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 80
$ws.Range("H80").Value = 54996
$ws.Range("J80").Value = 54996
$ws.Range("L80").Value = 54996
$ws.Range("N80").Value = -56992

# Row 83
$ws.Range("H83").Value = 54996
$ws.Range("J83").Value = 54996
$ws.Range("L83").Value = 164988
$ws.Range("N83").Value = -174972

# Row 107
$ws.Range("H107").Value = 40152
$ws.Range("J107").Value = 40152
$ws.Range("L107").Value = 40152
$ws.Range("N107").Value = -47832

# Row 109
$ws.Range("H109").Value = 38748.668
$ws.Range("J109").Value = 38748.668
$ws.Range("L109").Value = 38748.668
$ws.Range("N109").Value = -41522.668

# Row 111
$ws.Range("H111").Value = 48494
$ws.Range("J111").Value = 48494
$ws.Range("L111").Value = 48494
$ws.Range("N111").Value = -56674

# Row 114
$ws.Range("H114").Value = 44942
$ws.Range("J114").Value = 44942
$ws.Range("L114").Value = 44942
$ws.Range("N114").Value = -53620

# Row 118
$ws.Range("H118").Value = 49803
$ws.Range("J118").Value = 49803
$ws.Range("L118").Value = 49803
$ws.Range("N118").Value = -53117

# Row 119
$ws.Range("H119").Value = 50689.5
$ws.Range("J119").Value = 50689.5
$ws.Range("L119").Value = 50689.5
$ws.Range("N119").Value = -60365.5

# Row 121
$ws.Range("H121").Value = 44559.75
$ws.Range("J121").Value = 44559.75
$ws.Range("L121").Value = 44559.75
$ws.Range("N121").Value = -48053.75

# Row 123
$ws.Range("H123").Value = 46996
$ws.Range("J123").Value = 46996
$ws.Range("L123").Value = 46996
$ws.Range("N123").Value = -56796

# Row 125
$ws.Range("H125").Value = 44676.75
$ws.Range("J125").Value = 44676.75
$ws.Range("L125").Value = 44676.75
$ws.Range("N125").Value = -54516.75

# Row 134
$ws.Range("H134").Value = 52242.855
$ws.Range("J134").Value = 52242.855
$ws.Range("L134").Value = 52242.855
$ws.Range("N134").Value = -62382.855


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 119
$ws.Range("H119").Value = 41171.332
$ws.Range("J119").Value = 41171.332
$ws.Range("L119").Value = 41171.332
$ws.Range("N119").Value = -50847.332

# Row 120
$ws.Range("H120").Value = 48753
$ws.Range("J120").Value = 48753
$ws.Range("L120").Value = 48753
$ws.Range("N120").Value = -58429

# Row 137
$ws.Range("H137").Value = 40643.625
$ws.Range("J137").Value = 40643.625
$ws.Range("L137").Value = 40643.625
$ws.Range("N137").Value = -50843.625


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 45908.4
$ws.Range("J20").Value = 45908.4
$ws.Range("L20").Value = 45908.4
$ws.Range("N20").Value = -46380.4

# Row 30
$ws.Range("H30").Value = 45908.4
$ws.Range("J30").Value = 45908.4
$ws.Range("L30").Value = 45908.4
$ws.Range("N30").Value = -46090.4

# Row 100
$ws.Range("H100").Value = 29073.5
$ws.Range("J100").Value = 32431.334
$ws.Range("L100").Value = 32431.334
$ws.Range("N100").Value = -34595.334

# Row 116
$ws.Range("H116").Value = 48246.668
$ws.Range("J116").Value = 48246.668
$ws.Range("L116").Value = 48246.668
$ws.Range("N116").Value = -57424.668

# Row 118
$ws.Range("H118").Value = 40972.25
$ws.Range("J118").Value = 40972.25
$ws.Range("L118").Value = 40972.25
$ws.Range("N118").Value = -44286.25

# Row 119
$ws.Range("H119").Value = 48999.5
$ws.Range("J119").Value = 48999.5
$ws.Range("L119").Value = 48999.5
$ws.Range("N119").Value = -58675.5

# Row 128
$ws.Range("H128").Value = 45908.4
$ws.Range("J128").Value = 45908.4
$ws.Range("L128").Value = 45908.4
$ws.Range("N128").Value = -55868.4

# Row 141
$ws.Range("H141").Value = 12498.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 12498.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 12498.5
$ws.Range("N141").Value = -22858.5
$ws.Range("M141").ClearContents()


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 130
$ws.Range("H130").Value = 48522.75
$ws.Range("J130").Value = 48522.75
$ws.Range("L130").Value = 48522.75
$ws.Range("N130").Value = -58562.75

# Row 135
$ws.Range("H135").Value = 42179.8
$ws.Range("J135").Value = 42179.8
$ws.Range("L135").Value = 42179.8
$ws.Range("N135").Value = -52319.8


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# Row 114
$ws.Range("H114").Value = 38386
$ws.Range("J114").Value = 38386
$ws.Range("L114").Value = 38386
$ws.Range("N114").Value = -47064

# Row 116
$ws.Range("H116").Value = 50668
$ws.Range("J116").Value = 50668
$ws.Range("L116").Value = 50668
$ws.Range("N116").Value = -59846

# Row 121
$ws.Range("H121").Value = 43420
$ws.Range("J121").Value = 43420
$ws.Range("L121").Value = 43420
$ws.Range("N121").Value = -46914

# Row 127
$ws.Range("H127").Value = 41331.855
$ws.Range("I127").Value = 10650
$ws.Range("J127").Value = 46445.5
$ws.Range("K127").Value = 10650
$ws.Range("L127").Value = 46445.5
$ws.Range("N127").Value = -56365.5
$ws.Range("M127").Value = -5690

# Row 131
$ws.Range("H131").Value = 43322
$ws.Range("J131").Value = 43322
$ws.Range("L131").Value = 43322
$ws.Range("N131").Value = -53402


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 108
$ws.Range("H108").Value = 38528
$ws.Range("J108").Value = 38528
$ws.Range("L108").Value = 38528
$ws.Range("N108").Value = -46208

# Row 110
$ws.Range("H110").Value = 48644
$ws.Range("J110").Value = 48644
$ws.Range("L110").Value = 48644
$ws.Range("N110").Value = -56824

# Row 116
$ws.Range("H116").Value = 49672
$ws.Range("J116").Value = 49672
$ws.Range("L116").Value = 49672
$ws.Range("N116").Value = -58850

# Row 119
$ws.Range("H119").Value = 48693.5
$ws.Range("J119").Value = 48693.5
$ws.Range("L119").Value = 48693.5
$ws.Range("N119").Value = -58369.5

# Row 120
$ws.Range("H120").Value = 42037.332
$ws.Range("J120").Value = 42037.332
$ws.Range("L120").Value = 42037.332
$ws.Range("N120").Value = -51713.332

# Row 121
$ws.Range("H121").Value = 43246.668
$ws.Range("J121").Value = 43246.668
$ws.Range("L121").Value = 43246.668
$ws.Range("N121").Value = -46740.668

# Row 123
$ws.Range("H123").Value = 32510.2
$ws.Range("J123").Value = 31887.75
$ws.Range("L123").Value = 31887.75
$ws.Range("N123").Value = -41687.75

# Row 133
$ws.Range("H133").Value = 85510.5
$ws.Range("J133").Value = 85510.5
$ws.Range("L133").Value = 85510.5
$ws.Range("N133").Value = -95630.5

# Row 137
$ws.Range("H137").Value = 63999.5
$ws.Range("J137").Value = 63999.5
$ws.Range("L137").Value = 63999.5
$ws.Range("N137").Value = -74199.5

